# "Revert 'Revert "update project members"'" -- restore the two project
# members (Nabin Sharma, Daniel Maurath) that were previously dropped from
# the Sheet1 roster, re-adding their rows (with mailto hyperlinks on the
# "Preferred Contact Email" column) below the existing Pratik Mehta entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing header row (1) / Pratik Mehta row (2) keep their content as-is;
# just normalize the row height the same way the rest of the sheet uses.
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15

# --- Row 3: Nabin Sharma ---------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Nabin Sharma"
$ws.Cells.Item(3, 2).Value = "Providence, RI, USA"
$ws.Cells.Item(3, 3).Value = "EST (GMT – 0400 hrs)"
$ws.Cells.Item(3, 4).Value = "nabin.s.sharma@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:nabin.s.sharma@gmail.com")
$ws.Range("D3").Style = "Normal"
$ws.Cells.Item(3, 5).Value = "Data Science, Programming, DSP"
$ws.Cells.Item(3, 6).Value = "Software Developer"
$ws.Rows.Item(3).RowHeight = 15

# --- Row 4: Daniel Maurath --------------------------------------------------
$ws.Cells.Item(4, 1).Value = "Daniel Maurath"
$ws.Cells.Item(4, 2).Value = "San Francisco, CA"
$ws.Cells.Item(4, 3).Value = "PST (GMT - 0700 hrs)"
$ws.Cells.Item(4, 4).Value = "dmaurath@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:dmaurath@gmail.com")
$ws.Range("D4").Style = "Hyperlink"
$ws.Cells.Item(4, 5).Value = "Big Data, Predictive Analytics, Worforce Science"
$ws.Cells.Item(4, 6).Value = "Graduate Student in I/O Psychology"

# Match the saved selection/cursor position left behind in the sheet.
[void]$ws.Range("E9").Select()
